# Round the coordinate values in Q4/R4 to whole numbers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q4").Value = 613881
$ws.Range("R4").Value = 7034406

# Clear the Starttid (Z4) and Sluttid (AB4) time values, leaving the
# Startdatum (Y4) and Slutdatum (AA4) date values untouched.
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
